$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 70000
$ws.Range("D2").Value = 0.44461726443181
$ws.Range("E2").Value = 4.016195582332193
$ws.Range("F2").Value = 0.6898
$ws.Range("H2").Value = 4.16155928477495

$ws.Range("B3").Value = 70000
$ws.Range("D3").Value = 0.6004381621718412
$ws.Range("E3").Value = 3.972329726246276
$ws.Range("F3").Value = 1.1192
$ws.Range("H3").Value = 4.16155928477495

$ws.Range("B4").Value = 70000
$ws.Range("D4").Value = 0.7622662883774594
$ws.Range("E4").Value = 3.894095183543521
$ws.Range("F4").Value = 1.6471
$ws.Range("H4").Value = 4.16155928477495
